# Update: aliasing from single crystal reflection
#
# 1) Footer "updates automatically" date placeholders roll from 2021/2/25
#    to 2021/2/26 (slide master, every custom layout, and the notes master).
# 2) Slide 10 title "Issues" -> "In progress"
# 3) Slide 11 title "Issues, output" -> "Stretcher/compressor output"
# 4) Slide 2 outline bullet "Issues" -> "In progress"
# 5) Slide 7 red caption: "... (Si 444, 29.5 degree)" -> "... (Si 444, -29.5 degree)"

$p = $ppt.ActivePresentation

$oldDate = "2021/2/25"
$newDate = "2021/2/26"

# --- 1a. Slide master: find the date placeholder shape and update its text.
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $shape = $sm.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 1b. Every custom (slide) layout has its own copy of the date placeholder.
$layouts = $sm.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                    $shape.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# --- 1c. Notes master date placeholder only accepts edits through the
#     HeadersFooters façade in this host (the backing shape's TextFrame is
#     read-only here, and DateAndTime.Text reads back blank before the
#     first write even though the cached field text is $oldDate).
$nm = $p.NotesMaster
$nmDt = $nm.HeadersFooters.DateAndTime
$nmDt.Text = $newDate

# --- 2. Slide 10 title.
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(2).TextFrame.TextRange.Text = "In progress"

# --- 3. Slide 11 title.
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item(3).TextFrame.TextRange.Text = "Stretcher/compressor output"

# --- 4. Slide 2 outline bullet: replace just the "Issues" word in place so
#     the rest of the bulleted list is untouched.
$s2 = $p.Slides.Item(2)
$s2Body = $s2.Shapes.Item(2).TextFrame.TextRange
$full = $s2Body.Text
$needle = "Issues"
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $s2Body.Characters($idx + 1, $needle.Length).Text = "In progress"
}

# --- 5. Slide 7 red caption text box.
$s7 = $p.Slides.Item(7)
$s7Caption = $s7.Shapes.Item(4).TextFrame.TextRange
$s7Caption.Text = "100fs, 9481eV, one asymmetric reflection (Si 444, -29.5 degree)"
